$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sample identifiers in column A (A2:A25) to match expression matrix naming
$ws.Range("A2").Value = "0301C1"
$ws.Range("A3").Value = "0302C2"
$ws.Range("A4").Value = "0303C3"
$ws.Range("A5").Value = "0304X"
$ws.Range("A6").Value = "0305Y"
$ws.Range("A7").Value = "0306Z"
$ws.Range("A8").Value = "1307C1"
$ws.Range("A9").Value = "1308C2"
$ws.Range("A10").Value = "1309C3"
$ws.Range("A11").Value = "1310X"
$ws.Range("A12").Value = "1311Y"
$ws.Range("A13").Value = "1312Z"
$ws.Range("A14").Value = "2313C1"
$ws.Range("A15").Value = "2314C2"
$ws.Range("A16").Value = "2315C3"
$ws.Range("A17").Value = "2316X"
$ws.Range("A18").Value = "2317Y"
$ws.Range("A19").Value = "2318Z"
$ws.Range("A20").Value = "4325C1"
$ws.Range("A21").Value = "4326C2"
$ws.Range("A22").Value = "4327C3"
$ws.Range("A23").Value = "4328X"
$ws.Range("A24").Value = "4329Y"
$ws.Range("A25").Value = "4330Z"

# Update selected / visible cell state to reflect post-edit view
$ws.Range("A26").Select()
